# Atualização automática de preços de eletricidade
# Update row 2 (day 2025-...) of the Spot_PT sheet with the new hourly prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45926
$ws.Range("B2").Value = 106.48
$ws.Range("C2").Value = 104.52
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 98.70999999999999
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 104.52
$ws.Range("H2").Value = 107.96
$ws.Range("I2").Value = 123.37
$ws.Range("J2").Value = 123.28
$ws.Range("K2").Value = 108.8
$ws.Range("L2").Value = 88.59
$ws.Range("M2").Value = 34.99
$ws.Range("N2").Value = 25.76
$ws.Range("O2").Value = 29.2
$ws.Range("P2").Value = 21.56
$ws.Range("Q2").Value = 47.06
$ws.Range("R2").Value = 53
$ws.Range("S2").Value = 59.28
$ws.Range("T2").Value = 94.08
$ws.Range("U2").Value = 123.6
$ws.Range("V2").Value = 150
$ws.Range("W2").Value = 160
$ws.Range("X2").Value = 122.01
$ws.Range("Y2").Value = 113.92
$ws.Range("Z2").Value = 91.65000000000001
$ws.Range("AB2").Value = 136.48
$ws.Range("AD2").Value = 155
$ws.Range("AF2").Value = 117.96
$ws.Range("AG2").Value = "10h-17h"
